# Applies the horizontal.docx template tweak:
#   - nested "DescAndWeight" table row: trHeight 1933 -> 1846 twips
#   - nested "Price / DOH / Ratio" table: column widths 900/1800 -> 810/1890
#     (gridCol + matching tcW on the 2nd/3rd cell) and the 3rd cell's
#     right cell margin 202 -> 72 twips
#
# This runtime's Tables/Rows/Cells COM collections do not correctly
# resolve nested tables (a Table's .Tables/.Cell(...).Tables always
# loops back to the outermost table), so row-height/column-width
# properties can't be set through Row.Height / Cell.Width here.
# Instead we rebuild the body's WordprocessingML with the exact
# targeted substitutions and push it back with Range.InsertXML, which
# replaces the contents of the given Range with the supplied WordML -
# using the whole-document Range performs a clean whole-body swap.

$d = $word.ActiveDocument

$body = @'
<w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblpPr w:vertAnchor="page" w:horzAnchor="margin" w:tblpXSpec="center" w:tblpY="4494"/><w:tblOverlap w:val="never"/><w:tblW w:w="4896" w:type="dxa"/><w:tblBorders><w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tblBorders><w:tblLayout w:type="fixed"/><w:tblCellMar><w:left w:w="72" w:type="dxa"/><w:right w:w="15" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="0000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="4896"/></w:tblGrid><w:tr w:rsidR="00261B7D" w:rsidRPr="0093464C" w14:paraId="5422EF4F" w14:textId="77777777" w:rsidTr="000F5421"><w:trPr><w:cantSplit/><w:trHeight w:hRule="exact" w:val="3456"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4752" w:type="dxa"/><w:noWrap/><w:vAlign w:val="center"/></w:tcPr><w:tbl><w:tblPr><w:tblpPr w:leftFromText="187" w:rightFromText="187" w:vertAnchor="page" w:horzAnchor="margin" w:tblpXSpec="center" w:tblpY="93"/><w:tblOverlap w:val="never"/><w:tblW w:w="4608" w:type="dxa"/><w:jc w:val="center"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:left w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:right w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/></w:tblBorders><w:tblCellMar><w:top w:w="115" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="0000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="4608"/></w:tblGrid><w:tr w:rsidR="002879DC" w:rsidRPr="0093464C" w14:paraId="7F619777" w14:textId="77777777" w:rsidTr="00DA3BD9"><w:trPr><w:cantSplit/><w:trHeight w:hRule="exact" w:val="1933"/><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4608" w:type="dxa"/><w:tcBorders><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:tcMar><w:top w:w="43" w:type="dxa"/><w:left w:w="115" w:type="dxa"/><w:right w:w="115" w:type="dxa"/></w:tcMar></w:tcPr><w:p w14:paraId="797E5A78" w14:textId="5F039392" w:rsidR="005F7E68" w:rsidRPr="0093464C" w:rsidRDefault="005F7E68" w:rsidP="009D2E74"><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="00CC51E7"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>{{Label1.DescAndWeight}}</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p w14:paraId="5C372C8A" w14:textId="77777777" w:rsidR="0034354B" w:rsidRDefault="00200931" w:rsidP="009D2E74"><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:ind w:right="126"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>{{</w:t></w:r><w:r w:rsidR="00053E0C" w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Label1</w:t></w:r><w:r w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00053E0C" w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Lineage</w:t></w:r><w:r w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>}}</w:t></w:r></w:p><w:p w14:paraId="22C2A3DD" w14:textId="4D80AB44" w:rsidR="00E80EEC" w:rsidRPr="005A7DEA" w:rsidRDefault="002F4FA8" w:rsidP="009D2E74"><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:ind w:right="126"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="008F15C5"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="8"/><w:szCs w:val="8"/></w:rPr><w:t>{{Label1.ProductStrain}}</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblpPr w:leftFromText="187" w:rightFromText="187" w:vertAnchor="page" w:horzAnchor="margin" w:tblpX="90" w:tblpY="2452"/><w:tblOverlap w:val="never"/><w:tblW w:w="4590" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:left w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:right w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="FFFFFF" w:themeColor="background1"/></w:tblBorders><w:tblCellMar><w:right w:w="58" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="0000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="1890"/><w:gridCol w:w="900"/><w:gridCol w:w="1800"/></w:tblGrid><w:tr w:rsidR="00BE6EA5" w:rsidRPr="0093464C" w14:paraId="003902D7" w14:textId="77777777" w:rsidTr="00DA3BD9"><w:trPr><w:cantSplit/><w:trHeight w:hRule="exact" w:val="897"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="1890" w:type="dxa"/><w:tcBorders><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:noWrap/><w:tcMar><w:top w:w="72" w:type="dxa"/><w:left w:w="115" w:type="dxa"/><w:right w:w="14" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="3F36A83C" w14:textId="77777777" w:rsidR="00BE6EA5" w:rsidRPr="0093464C" w:rsidRDefault="00BE6EA5" w:rsidP="00BE6EA5"><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:ind w:right="126"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="0093464C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/></w:rPr><w:t>{{Label1.Price}}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="900" w:type="dxa"/><w:tcBorders><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:noWrap/><w:tcMar><w:top w:w="58" w:type="dxa"/><w:left w:w="115" w:type="dxa"/></w:tcMar><w:vAlign w:val="bottom"/></w:tcPr><w:p w14:paraId="05C4CCCF" w14:textId="3C851647" w:rsidR="00BE6EA5" w:rsidRPr="000D6574" w:rsidRDefault="00BE6EA5" w:rsidP="0032333D"><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:ind w:right="126"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="11"/><w:szCs w:val="11"/></w:rPr></w:pPr><w:r w:rsidRPr="000D6574"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="11"/><w:szCs w:val="11"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/></w:rPr><w:t>{{Label1.DOH}}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1800" w:type="dxa"/><w:tcBorders><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:noWrap/><w:tcMar><w:top w:w="72" w:type="dxa"/><w:left w:w="72" w:type="dxa"/><w:right w:w="202" w:type="dxa"/></w:tcMar><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="69F1E708" w14:textId="77777777" w:rsidR="00BE6EA5" w:rsidRPr="004F783E" w:rsidRDefault="00BE6EA5" w:rsidP="00A62B42"><w:pPr><w:spacing w:line="264" w:lineRule="auto"/><w:ind w:right="130"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="004F783E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>{{Label1.Ratio_o</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>r</w:t></w:r><w:r w:rsidRPr="004F783E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>_THC_CBD }}</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p w14:paraId="3C4EA489" w14:textId="511523CF" w:rsidR="007B40F4" w:rsidRPr="00FB27ED" w:rsidRDefault="007B40F4" w:rsidP="009D2E74"><w:pPr><w:spacing w:line="216" w:lineRule="auto"/><w:ind w:right="126"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="2"/><w:szCs w:val="2"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p w14:paraId="4502C4E1" w14:textId="75A13736" w:rsidR="009F11B8" w:rsidRPr="0093464C" w:rsidRDefault="009F11B8" w:rsidP="00666E08"><w:pPr><w:tabs><w:tab w:val="left" w:pos="10937"/></w:tabs><w:ind w:left="126" w:right="126"/><w:rPr><w:vanish/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p><w:sectPr w:rsidR="009F11B8" w:rsidRPr="0093464C" w:rsidSect="007378E4"><w:type w:val="continuous"/><w:pgSz w:w="15840" w:h="12240" w:orient="landscape"/><w:pgMar w:top="432" w:right="806" w:bottom="360" w:left="360" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/></w:sectPr>
'@

# --- 1) nested table row height: 1933 -> 1846 (trHeight, hRule="exact") ---
$old1 = '<w:trHeight w:hRule="exact" w:val="1933"/>'
$new1 = '<w:trHeight w:hRule="exact" w:val="1846"/>'
if ($body.IndexOf($old1) -lt 0) { throw "pattern not found: $old1" }
$body = $body.Replace($old1, $new1)

# --- 2) tblGrid column widths: 900/1800 -> 810/1890 ---
$old2 = '<w:gridCol w:w="1890"/><w:gridCol w:w="900"/><w:gridCol w:w="1800"/>'
$new2 = '<w:gridCol w:w="1890"/><w:gridCol w:w="810"/><w:gridCol w:w="1890"/>'
if ($body.IndexOf($old2) -lt 0) { throw "pattern not found: $old2" }
$body = $body.Replace($old2, $new2)

# --- 3) 2nd column cell width: 900 -> 810 ---
$old3 = '<w:tcW w:w="900" w:type="dxa"/>'
$new3 = '<w:tcW w:w="810" w:type="dxa"/>'
if ($body.IndexOf($old3) -lt 0) { throw "pattern not found: $old3" }
$body = $body.Replace($old3, $new3)

# --- 4) 3rd column cell width: 1800 -> 1890 ---
$old4 = '<w:tcW w:w="1800" w:type="dxa"/>'
$new4 = '<w:tcW w:w="1890" w:type="dxa"/>'
if ($body.IndexOf($old4) -lt 0) { throw "pattern not found: $old4" }
$body = $body.Replace($old4, $new4)

# --- 5) 3rd column cell right margin: 202 -> 72 ---
$old5 = '<w:right w:w="202" w:type="dxa"/>'
$new5 = '<w:right w:w="72" w:type="dxa"/>'
if ($body.IndexOf($old5) -lt 0) { throw "pattern not found: $old5" }
$body = $body.Replace($old5, $new5)

$d.Content.InsertXML($body)

Write-Host "Applied horizontal.docx sizing tweaks."
